{"js": "// Remove the empty \"FirstParagraph\"-styled paragraph (containing only a\n// manual line break) that immediately follows the \"13 Esophagogastrectomy\"\n// Heading2 paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet target = null;\n\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"13 Esophagogastrectomy\") {\n    const next = items[i + 1];\n    // Guard: only remove it if it is the expected blank \"line break only\"\n    // paragraph right after the heading (style \"First Paragraph\", text is\n    // just the vertical-tab break char, or empty).\n    if (next && next.style === \"First Paragraph\" && next.text.replace(/\\v/g, \"\").trim() === \"\") {\n      target = next;\n    }\n    break;\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the empty \"First Paragraph\"-styled paragraph (containing only a\n# manual line break) that immediately follows the \"13 Esophagogastrectomy\"\n# Heading 2 paragraph.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($t -eq \"13 Esophagogastrectomy\") {\n        if ($i + 1 -le $count) {\n            $next = $d.Paragraphs.Item($i + 1)\n            $nextStyleName = $next.Style.NameLocal\n            $nextText = $next.Range.Text.TrimEnd([char]13, [char]7)\n\n            # Guard: only remove it if it is the expected blank \"line break\n            # only\" paragraph right after the heading (style \"First\n            # Paragraph\", text is just the vertical-tab break char, or empty).\n            if ($nextStyleName -eq \"First Paragraph\" -and $nextText.Trim([char]11) -eq \"\") {\n                $next.Range.Delete()\n            }\n        }\n        break\n    }\n}\n"}
